$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

# Title (appears twice: H1 heading and bold summary line near the end)
Replace-Text "Play Feng Fu for Free - Discover Chinese Symbolism and Lucky Animals" "Play Feng Fu Slot Game for Free - Soak in Chinese Symbolism"

# "What we like" bullet list
Replace-Text "Chinese-themed design and lucky symbols enhance gameplay experience" "Chinese symbolism related to luck and wealth"
Replace-Text "Use of Chinese ideograms adds an extra layer of depth and meaning" "Lucky animals and Chinese zodiac theme"
Replace-Text "Traditional slot game mechanics are ideal for experienced slot players" "Use of Chinese ideograms adds depth to the game"
Replace-Text "Bonus function adds another opportunity to win big" "Challenging gameplay for experienced slot players"

# "What we don't like" bullet list
Replace-Text "Chinese ideograms can be confusing for new players" "Difficulty for new players unfamiliar with Chinese ideograms"
Replace-Text "Limited pay lines may not offer enough variety for some slot enthusiasts" "Limited accessibility for non-Chinese culture enthusiasts"

# Closing italic summary line
Replace-Text "Explore traditional gameplay mechanics and try your luck with Chinese-themed symbols. Play Feng Fu for free and win big with bonus functions." "Play Feng Fu for free and explore Chinese symbolism related to luck and wealth."
